# Update database and change read_price algorithm
#
# The workbook tracks five twelve-month periods ending 1396/12 .. 1400/12 in
# columns E..I. This edit rolls the window forward by one year: the oldest
# period (1396/12) is dropped, the remaining four periods shift one column to
# the left (E<-F, F<-G, G<-H, H<-I), and a new period (1401/12) is appended in
# column I with freshly read figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column headers (row 8 and row 24) -------------------------------------
# Shift the "twelve months ending ..." labels one column to the left and add
# the new 1401/12 label in column I. This makes the old 1396/12 label
# unreferenced so it drops out of the shared-string table, and introduces the
# new 1401/12 label.
$ws.Range("E8").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E24").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F24").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G24").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H24").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I24").Value = "دوازده ماهه منتهی به 1401/12"

# --- Data rows ---------------------------------------------------------------
# Helper: shift row r's E:I values one column left, and place $newI in I.
function Shift-Row($r, $newI) {
    $e = $ws.Range("F$r").Value2
    $f = $ws.Range("G$r").Value2
    $g = $ws.Range("H$r").Value2
    $h = $ws.Range("I$r").Value2
    $ws.Range("E$r").Value = $e
    $ws.Range("F$r").Value = $f
    $ws.Range("G$r").Value = $g
    $ws.Range("H$r").Value = $h
    $ws.Range("I$r").Value = $newI
}

Shift-Row 10 0
Shift-Row 11 0
Shift-Row 12 0
Shift-Row 13 0
Shift-Row 14 33195
Shift-Row 15 15447
Shift-Row 16 102942
Shift-Row 17 915825
Shift-Row 18 0
Shift-Row 19 8882829
Shift-Row 20 9950238
Shift-Row 26 5704
Shift-Row 27 235
